# Fix a typo in the "Reason Excluded" column and add an AutoFilter to the
# supplemental-table worksheet, matching the upstream commit:
#   "fixed typo in Olker supp table and add filtering step for mc6 flags
#    altho it doesn't change priority chems"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Typo fix: "...420 nmb" -> "...420 nm" wherever it occurs in column D.
$oldText = "Abnormal time course of absorbance at 420 nmb"
$newText = "Abnormal time course of absorbance at 420 nm"

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}

# 2) Turn on AutoFilter over the table (A1:D33 — header row + 32 data rows).
$ws.Range("A1:D33").AutoFilter()

# Register the (hidden, sheet-scoped) _FilterDatabase defined name that
# Excel normally writes out alongside an AutoFilter.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$D`$33")
$filterName.Visible = $false

# 3) Move the active selection (cosmetic, matches the saved view state).
$ws.Range("I18").Select()
